$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bank the existing "Hyperlink" cell formatting (font + border) from C3 / C5 / C6 onto scratch
# cells so it can be restored after the hyperlinks are recreated below (Hyperlinks.Add() resets
# cell formatting to a fresh Hyperlink style that drops the border the author had added on top).
$ws.Range("C3").Copy()
$ws.Range("Z1").PasteSpecial(-4122)
$ws.Range("C5").Copy()
$ws.Range("Z2").PasteSpecial(-4122)
$ws.Range("C6").Copy()
$ws.Range("Z3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 6 (TC001 / positive): email changes to salsabilarhdsy@gmail.com ---
$ws.Range("C6").Value = "salsabilarhdsy@gmail.com"

# --- Row 3 (TC003 / negative): email changes to salsablsy@gmail.com, password to Salsa123! ---
$ws.Range("C3").Value = "salsablsy@gmail.com"
$ws.Range("D3").Value = "Salsa123!"

# --- Row 4 (TC004 / negative): email changes to salsablsy, password to Salsa123! ---
$ws.Range("C4").Value = "salsablsy"
$ws.Range("D4").Value = "Salsa123!"

# D3/D4 drop the "vertical-top" alignment override and match D5's plain bordered style. Copy the
# format from D5 (already-existing style in the workbook) instead of assigning a brand-new style
# object, so the cellXfs table doesn't grow for this part of the edit.
$ws.Range("D5").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("D4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 5 (TC005 / negative): email changes to salsabilarhdsy@gmail.com, password to Salsa123* ---
$ws.Range("C5").Value = "salsabilarhdsy@gmail.com"
$ws.Range("D5").Value = "Salsa123*"

# --- Hyperlinks: re-point mailto targets at the new e-mail addresses, keep C6/C3/C5 order ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C6"), "mailto:salsabilarhdsy@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:salsablsy@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:salsabilarhdsy@gmail.com")

# restore the banked formatting clobbered by Hyperlinks.Add()
$ws.Range("Z1").Copy()
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("Z2").Copy()
$ws.Range("C5").PasteSpecial(-4122)
$ws.Range("Z3").Copy()
$ws.Range("C6").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("Z1:Z3").Clear() | Out-Null

# --- Column C widened to fit the longer e-mail addresses ---
$ws.Columns.Item(3).ColumnWidth = 24.6

# --- Selection moved down one row ---
$ws.Range("E12").Select() | Out-Null
